$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M6").Value = 5.76
$wsGrupo.Range("M24").Value = "7 de 22"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F6").Value = 5.76
$wsMensual.Range("F24").Value = 30732.96

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D16").Value = 27195.98
$wsCumpl.Range("E16").Value = 11560.56
$wsCumpl.Range("F16").Value = 0.7017133108373451

$wsCumpl.Range("D19").Value = 30732.96
$wsCumpl.Range("E19").Value = 27490.04386304603
$wsCumpl.Range("F19").Value = 0.5278490967640733
